# "adjustments to table and explanation"
#
# Loop Table 2: the D/E/F columns were only ever showing the single
# character contributed by the current loop iteration; they should show
# the *running* (cumulative) value of s1/s2 instead. Also refresh the
# long walkthrough note in G14 to match, shrink that row now that the
# note is a little shorter, and update the saved selection/scroll spot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loop Table 2")
$ws.Activate()

# --- table body: cumulative s1/s2 strings instead of single letters ---
$ws.Range("D16").Value = "me"
$ws.Range("F16").Value = "me"

$ws.Range("E17").Value = "pu"
$ws.Range("F17").Value = "pu"

$ws.Range("D18").Value = "meo"
$ws.Range("F18").Value = "meo"

$ws.Range("E19").Value = "pur"
$ws.Range("F19").Value = "pur"

$ws.Range("D20").Value = "meow"
$ws.Range("F20").Value = "meow"

$ws.Range("E21").Value = "purr"
$ws.Range("F21").Value = "purr"

# --- updated walkthrough explanation in G14 ---
$ws.Range("G14").Value = "index starts at 0, and the loop will continue as long as index is less than 8 (alpha.length);  if index divided by 2 results in a remainder of zero (meaning index is even), s1 is equal to s1 plus the string character that is the value of index, 1; if it's an odd number (modulo is not 0), s2 is equal to the value of s2 plus the string character that is the value of index, 1; at the end of every loop we'll add one to index; when the condition is no longer true (when index > alpha.length), and then it will display s1 + s2 in console --> ""meow"" + ""purr"" = meowpurr"

# --- row 14 is shorter now that the note reads more tersely ---
$ws.Rows.Item(14).RowHeight = 96

# --- view: scrolled up a couple of rows, selection moved up one row ---
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("G16").Select()

# --- application window nudged slightly to match the saved layout ---
$win.Left = 8520
$win.Top = 3240
$win.Width = 22520
$win.Height = 12740
